# Update "want to go" counts (column F) on several rows across sheets to
# reflect refreshed data generated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F19").Value = 10795
$ws1.Range("F20").Value = 6206
$ws1.Range("F22").Value = 13
$ws1.Range("F23").Value = 398
$ws1.Range("F27").Value = 858
$ws1.Range("F28").Value = 31
$ws1.Range("F29").Value = 202
$ws1.Range("F32").Value = 51
$ws1.Range("F37").Value = 252
$ws1.Range("F38").Value = 263
$ws1.Range("F39").Value = 4890
$ws1.Range("F41").Value = 1168

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F12").Value = 117

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 8891
$ws3.Range("F3").Value = 453

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 8891
$ws4.Range("F3").Value = 453
$ws4.Range("F22").Value = 10795
$ws4.Range("F25").Value = 13
$ws4.Range("F26").Value = 398
$ws4.Range("F30").Value = 858
$ws4.Range("F31").Value = 31
$ws4.Range("F34").Value = 51
$ws4.Range("F40").Value = 263
$ws4.Range("F41").Value = 4890
$ws4.Range("F43").Value = 1168
